$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 63: set the date in column A and the activity description in column B
$ws.Range("A63").Value = Get-Date -Year 2013 -Month 4 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("B63").Value = "began with OpenCL chapter"

# Reflect the new selection on the sheet (user clicked/entered on B63)
$ws.Range("B63").Select()
